# Refresh the coin table: each row's Price (D) and Volume(1h) (E)
# columns get the latest scrape's figures. Rows 50/51 additionally
# swap which coin (Algorand vs RocketPoolETH) occupies that rank, so
# their Coin (B) and Link (C) cells are rewritten too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Prefix with a literal apostrophe -- Excel's "force text" marker --
    # so numeric-looking values (e.g. '306.10') are stored verbatim as
    # text instead of being coerced into a Double (which would also
    # silently drop the trailing zero).
    $ws.Range($cellRef).Value = "'" + $text
}

$ws.Range('D2').Value = '46.110.29'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '2.596.02'
$ws.Range('E3').Value = '  +8.75%  '
$ws.Range('E4').Value = '  -0.06%  '
Set-TextCell 'D5' '306.10'
$ws.Range('E5').Value = '  +1.45%  '
Set-TextCell 'D6' '99.83'
$ws.Range('E6').Value = '  +0.55%  '
Set-TextCell 'D7' '0.599'
$ws.Range('E7').Value = '  +5.28%  '
$ws.Range('E8').Value = '  -0.03%  '
Set-TextCell 'D9' '0.571'
$ws.Range('E9').Value = '  +11.62%  '
$ws.Range('E10').Value = '  +10.76%  '
Set-TextCell 'D11' '0.0833'
$ws.Range('E11').Value = '  +5.13%  '
Set-TextCell 'D12' '8.06'
$ws.Range('E12').Value = '  +12.57%  '
$ws.Range('D13').Value = '2.987.37'
$ws.Range('E13').Value = '  +8.43%  '
$ws.Range('E14').Value = '  +1.45%  '
$ws.Range('D15').Value = '2.606.46'
$ws.Range('E15').Value = '  +9.25%  '
Set-TextCell 'D16' '0.894'
$ws.Range('E16').Value = '  +8.54%  '
Set-TextCell 'D17' '14.82'
$ws.Range('E17').Value = '  +7.72%  '
$ws.Range('D18').Value = '46.224.65'
$ws.Range('E18').Value = '  +0.66%  '
Set-TextCell 'D19' '13.19'
$ws.Range('E19').Value = '  +3.11%  '
Set-TextCell 'D20' '0.0000100'
$ws.Range('E20').Value = '  +4.88%  '
$ws.Range('E21').Value = '  +9.55%  '
Set-TextCell 'D22' '70.73'
$ws.Range('E22').Value = '  +5.32%  '
Set-TextCell 'D23' '253.87'
$ws.Range('E23').Value = '  +3.78%  '
Set-TextCell 'D24' '2.98'
$ws.Range('E24').Value = '  +5.71%  '
Set-TextCell 'D25' '2.23'
$ws.Range('E25').Value = '  +15.02%  '
Set-TextCell 'D26' '27.87'
$ws.Range('E26').Value = '  +31.82%  '
Set-TextCell 'D27' '1.00'
$ws.Range('E27').Value = '  +0.01%  '
Set-TextCell 'D28' '10.44'
$ws.Range('E28').Value = '  +6.57%  '
Set-TextCell 'D29' '39.78'
$ws.Range('E29').Value = '  +0.13%  '
Set-TextCell 'D30' '2.25'
$ws.Range('E30').Value = '  +2.21%  '
Set-TextCell 'D31' '6.10'
$ws.Range('E31').Value = '  +9.70%  '
$ws.Range('E32').Value = '  -3.46%  '
Set-TextCell 'D33' '2.31'
$ws.Range('E33').Value = '  +18.47%  '
Set-TextCell 'D34' '2.90'
$ws.Range('E34').Value = '  +3.15%  '
Set-TextCell 'D35' '152.16'
$ws.Range('E35').Value = '  +3.64%  '
$ws.Range('E36').Value = '  +6.98%  '
$ws.Range('E37').Value = '  +2.88%  '
$ws.Range('E38').Value = '  +4.59%  '
Set-TextCell 'D39' '4.18'
$ws.Range('E39').Value = '  +6.35%  '
Set-TextCell 'D40' '15.50'
$ws.Range('E40').Value = '  +3.36%  '
Set-TextCell 'D41' '3.59'
$ws.Range('E41').Value = '  +9.63%  '
Set-TextCell 'D42' '0.0321'
$ws.Range('E42').Value = '  +6.65%  '
$ws.Range('D43').Value = '2.044.27'
$ws.Range('E43').Value = '  +5.13%  '
Set-TextCell 'D44' '19.02'
$ws.Range('E44').Value = '  +33.83%  '
Set-TextCell 'D45' '0.998'
$ws.Range('E45').Value = '  -0.05%  '
Set-TextCell 'D46' '90.83'
$ws.Range('E46').Value = '  -0.93%  '
Set-TextCell 'D47' '9.11'
$ws.Range('E47').Value = '  +7.45%  '
Set-TextCell 'D48' '109.20'
$ws.Range('E49').Value = '  -0.49%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 'D50' '0.199'
$ws.Range('E50').Value = '  +6.75%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.846.29'
$ws.Range('E51').Value = '  +8.40%  '
